$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25; this shifts the existing rows 25-75
# down to 26-76 (carrying their formatting/styles along), and grows the
# sheet dimension from R75 to R76 automatically.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly data record.
$ws.Range("A25").Value = 11
$ws.Range("B25").Value = "Vega Monumental Concepción"
$ws.Range("C25").Value = "Bíobío"
$ws.Range("D25").Value = 44519
$ws.Range("E25").Value = 8
$ws.Range("F25").Value = 100112032
$ws.Range("G25").Value = "Zapallo italiano"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 350
$ws.Range("K25").Value = 5500
$ws.Range("L25").Value = 6000
$ws.Range("M25").Value = 5786
$ws.Range("N25").Value = "$/caja 60 unidades"
$ws.Range("O25").Value = "Región de Arica y Parinacota"
$ws.Range("P25").Value = 96
$ws.Range("Q25").Value = 60
$ws.Range("R25").Value = "Hortaliza"
